$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.818.74"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "2.541.89"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.571"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.11%  "
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0805"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.24%  "
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").Value = "2.929.96"
$ws.Range("E14").Value = "  -1.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.77"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.28%  "
$ws.Range("D16").Value = "2.572.08"
$ws.Range("E16").Value = "  +2.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.815"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.49%  "
$ws.Range("D18").Value = "42.799.82"
$ws.Range("E18").Value = "  -0.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.15%  "
$ws.Range("D20").Value = "0.0₃0952"
$ws.Range("E20").Value = "  -1.12%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "244.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.07%  "
$ws.Range("E24").Value = "  -2.02%  "
$ws.Range("E25").Value = "  -1.75%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.12"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.16%  "
$ws.Range("E28").Value = "  -5.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "155.88"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0793"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.62"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.49%  "
$ws.Range("E35").Value = "  -5.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.36"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.19"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.112"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("E40").Value = "  +9.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.33"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("B43").Value = "NEARProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.24%  "
$ws.Range("E44").Value = "  -1.45%  "
$ws.Range("D45").Value = "1.977.68"
$ws.Range("E45").Value = "  -1.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("D47").Value = "2.782.10"
$ws.Range("E47").Value = "  -1.08%  "
$ws.Range("B48").Value = "SEI"
$ws.Range("C48").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.865"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +10.01%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.193"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "80.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.60%  "
